$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("C2").Value = 0.05088209513382935
$ws.Range("D2").Value = 0.09244243276017983
$ws.Range("E2").Value = 0.4343060183838787
$ws.Range("F2").Value = 2.285403960572268
$ws.Range("G2").Value = 2.564910388463545
$ws.Range("H2").Value = 1.534877293481372
$ws.Range("N2").Value = 2.873581687696344
# Row 3
$ws.Range("C3").Value = 0.04499111415316293
$ws.Range("D3").Value = 0.08099312179201945
$ws.Range("E3").Value = 0.3774527995275037
$ws.Range("F3").Value = 2.053156900224707
$ws.Range("G3").Value = 2.275893840270271
$ws.Range("H3").Value = 1.407448340503549
$ws.Range("N3").Value = 2.562605684679454
# Row 4
$ws.Range("C4").Value = 0.04140688078129529
$ws.Range("D4").Value = 0.07402400667390907
$ws.Range("E4").Value = 0.3428124898187548
$ws.Range("F4").Value = 1.91231671163419
$ws.Range("G4").Value = 2.100249222090838
$ws.Range("H4").Value = 1.330500229916026
$ws.Range("N4").Value = 2.371325805375761
# Row 5
$ws.Range("C5").Value = 0.03995402274580329
$ws.Range("D5").Value = 0.07119808863397736
$ws.Range("E5").Value = 0.3287567351206491
$ws.Range("F5").Value = 1.855343698008141
$ws.Range("G5").Value = 2.02910102924011
$ws.Range("H5").Value = 1.299455806609899
$ws.Range("N5").Value = 2.293303068605894
# Row 6
$ws.Range("C6").Value = 0.03971322982907566
$ws.Range("D6").Value = 0.07072966030489169
$ws.Range("E6").Value = 0.3264262411716459
$ws.Range("F6").Value = 1.845908133199856
$ws.Range("G6").Value = 2.017312013858941
$ws.Range("H6").Value = 1.29431941459319
$ws.Range("N6").Value = 2.280343261403857
# Row 7
$ws.Range("C7").Value = 0.04138725633737295
$ws.Range("D7").Value = 0.07398584006804754
$ws.Range("E7").Value = 0.3426226931926237
$ws.Range("F7").Value = 1.911546679905854
$ws.Range("G7").Value = 2.099287993742394
$ws.Range("H7").Value = 1.33008030506619
$ws.Range("N7").Value = 2.370273851392596
# Row 8
$ws.Range("C8").Value = 0.04884375063495838
$ws.Range("D8").Value = 0.08848133144674364
$ws.Range("E8").Value = 0.4146432302580223
$ws.Range("F8").Value = 2.204946202125313
$ws.Range("G8").Value = 2.464864458016393
$ws.Range("H8").Value = 1.490663432560154
$ws.Range("N8").Value = 2.766433886209825
# Row 9
$ws.Range("C9").Value = 0.06375108278101038
$ws.Range("D9").Value = 0.117447993817521
$ws.Range("E9").Value = 0.5583251547865871
$ws.Range("F9").Value = 2.795323943595179
$ws.Range("G9").Value = 3.197443359351496
$ws.Range("H9").Value = 1.816439590657581
$ws.Range("N9").Value = 3.540180268007646
# Row 10
$ws.Range("C10").Value = 0.07491319404151398
$ws.Range("D10").Value = 0.1391449197433303
$ws.Range("E10").Value = 0.6658540739309728
$ws.Range("F10").Value = 3.239773343520312
$ws.Range("G10").Value = 3.747132436533263
$ws.Range("H10").Value = 2.063304805149698
$ws.Range("N10").Value = 4.10623028343673
# Row 11
$ws.Range("C11").Value = 0.08004481728906399
$ws.Range("D11").Value = 0.1491249759874904
$ws.Range("E11").Value = 0.715306283212044
$ws.Range("F11").Value = 3.444639034242812
$ws.Range("G11").Value = 4.00012231043263
$ws.Range("H11").Value = 2.177445210765711
$ws.Range("N11").Value = 4.363110593465422
# Row 12
$ws.Range("C12").Value = 0.08199645954792345
$ws.Range("D12").Value = 0.1529216098636255
$ws.Range("E12").Value = 0.73411879543319
$ws.Range("F12").Value = 3.522631409936139
$ws.Range("G12").Value = 4.096380848005879
$ws.Range("H12").Value = 2.220948711345443
$ws.Range("N12").Value = 4.460285735713398
# Row 13
$ws.Range("C13").Value = 0.08157575385018845
$ws.Range("D13").Value = 0.1521031369487673
$ws.Range("E13").Value = 0.7300632132612037
$ws.Range("F13").Value = 3.505815472973268
$ws.Range("G13").Value = 4.075628964692669
$ws.Range("H13").Value = 2.211566688482378
$ws.Range("N13").Value = 4.439361943450422
# Row 14
$ws.Range("C14").Value = 0.08020520772986117
$ws.Range("D14").Value = 0.1494369692640589
$ws.Range("E14").Value = 0.7168522236950281
$ws.Range("F14").Value = 3.451047047003783
$ws.Range("G14").Value = 4.008032199225568
$ws.Range("H14").Value = 2.181018542013135
$ws.Range("N14").Value = 4.371107314139294
# Row 15
$ws.Range("C15").Value = 0.07936682344265478
$ws.Range("D15").Value = 0.1478061790211314
$ws.Range("E15").Value = 0.7087715839542454
$ws.Range("F15").Value = 3.417554602524831
$ws.Range("G15").Value = 3.966687771745626
$ws.Range("H15").Value = 2.162344012834808
$ws.Range("N15").Value = 4.329286057409945
# Row 16
$ws.Range("C16").Value = 0.07457897080199416
$ws.Range("D16").Value = 0.1384950473262165
$ws.Range("E16").Value = 0.6626338051536322
$ws.Range("F16").Value = 3.226441281925361
$ws.Range("G16").Value = 3.73066093909091
$ws.Range("H16").Value = 2.055883877073995
$ws.Range("N16").Value = 4.089429168003562
# Row 17
$ws.Range("C17").Value = 0.07165604684574589
$ws.Range("D17").Value = 0.1328122544846622
$ws.Range("E17").Value = 0.6344734356464699
$ws.Range("F17").Value = 3.109906547216241
$ws.Range("G17").Value = 3.586641909769298
$ws.Range("H17").Value = 1.99105681049997
$ws.Range("N17").Value = 3.94211849063862
# Row 18
$ws.Range("C18").Value = 0.06997989321288856
$ws.Range("D18").Value = 0.1295539192348656
$ws.Range("E18").Value = 0.618326308666667
$ws.Range("F18").Value = 3.0431296698263
$ws.Range("G18").Value = 3.504079926245083
$ws.Range("H18").Value = 1.953942185061351
$ws.Range("N18").Value = 3.857331695637129
# Row 19
$ws.Range("C19").Value = 0.06941322161689811
$ws.Range("D19").Value = 0.1284524172936017
$ws.Range("E19").Value = 0.6128674879012692
$ws.Range("F19").Value = 3.020562399879111
$ws.Range("G19").Value = 3.476171866273489
$ws.Range("H19").Value = 1.941404895010351
$ws.Range("N19").Value = 3.828614786364199
# Row 20
$ws.Range("C20").Value = 0.07196667133618462
$ws.Range("D20").Value = 0.1334161246630856
$ws.Range("E20").Value = 0.6374659200924242
$ws.Range("F20").Value = 3.122285708113651
$ws.Range("G20").Value = 3.601944379440397
$ws.Range("H20").Value = 1.99793982481799
$ws.Range("N20").Value = 3.957806003277312
# Row 21
$ws.Range("C21").Value = 0.08060753682887878
$ws.Range("D21").Value = 0.1502196021049826
$ws.Range("E21").Value = 0.7207302043796346
$ws.Range("F21").Value = 3.467122389459689
$ws.Range("G21").Value = 4.027874326515359
$ws.Range("H21").Value = 2.189983512906224
$ws.Range("N21").Value = 4.391158149571254
# Row 22
$ws.Range("C22").Value = 0.08630417082449071
$ws.Range("D22").Value = 0.1613038644014182
$ws.Range("E22").Value = 0.7756540621339951
$ws.Range("F22").Value = 3.694921204906962
$ws.Range("G22").Value = 4.308923702799916
$ws.Range("H22").Value = 2.317141151680005
$ws.Range("N22").Value = 4.673791817957863
# Row 23
$ws.Range("C23").Value = 0.08325903397826551
$ws.Range("D23").Value = 0.1553780927345656
$ws.Range("E23").Value = 0.7462908648888913
$ws.Range("F23").Value = 3.573108878531798
$ws.Range("G23").Value = 4.158665295903518
$ws.Range("H23").Value = 2.249118535815285
$ws.Range("N23").Value = 4.523002190001307
# Row 24
$ws.Range("C24").Value = 0.0718262247396666
$ws.Range("D24").Value = 0.1331430877325772
$ws.Range("E24").Value = 0.6361128857111993
$ws.Range("F24").Value = 3.116688404811185
$ws.Range("G24").Value = 3.595025398961184
$ws.Range("H24").Value = 1.994827531494195
$ws.Range("N24").Value = 3.950713976768498
# Row 25
$ws.Range("C25").Value = 0.05968409346765213
$ws.Range("D25").Value = 0.1095451899369664
$ws.Range("E25").Value = 0.5191453392989445
$ws.Range("F25").Value = 2.633846868522568
$ws.Range("G25").Value = 2.997389641848088
$ws.Range("H25").Value = 1.72705449418288
$ws.Range("N25").Value = 3.331249627311138
